$d = $word.ActiveDocument

# --- Locate the two paragraphs we need to touch, by content, so the
# script is resilient to exact paragraph numbering. ---
$para1 = $null   # "User1 sends out to server ..."
$para2 = $null   # "If user2 accepts ..."
$para3 = $null   # empty paragraph right after para2
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($para1 -eq $null -and $t -like "*User1 sends out*") {
        $para1 = $i
    }
    if ($para2 -eq $null -and $t -like "*If user2 accepts*") {
        $para2 = $i
        $para3 = $i + 1
    }
}

# --- 1) "User1 sends out ..." paragraph: split the trailing run so the
# "... and waits for response (01 ok, 02 reject)" tail is removed and the
# remaining text is split into two runs. Only the runs change; the
# paragraph's own pPr is left completely alone. ---
$rng1 = $d.Paragraphs($para1).Range
$xml1 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:right="90"/></w:pPr><w:r><w:t>User1 sends out to server (con,user2</w:t></w:r><w:r><w:t>uid</w:t></w:r><w:r><w:t xml:space="preserve">) and server sends </w:t></w:r><w:r><w:t xml:space="preserve">(req,user1uid) </w:t></w:r><w:r><w:t>to user2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1)

# --- 2) Replace the "If user2 accepts ..." paragraph and the empty
# paragraph that follows it with: the reworked "accepts" paragraph (no
# longer containing the bookmark), a new blank paragraph, a new "If
# user2 rejects ..." bullet (carrying the relocated bookmark), and the
# trailing blank paragraph whose indND changes from left=360/firstLine=0
# to just right=90. Re-resolve the paragraph objects/range fresh in case
# paragraph numbering shifted because of the first edit. ---
$para2 = $null
$para3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($para2 -eq $null -and $t -like "*If user2 accepts*") {
        $para2 = $i
        $para3 = $i + 1
    }
}
$rangeStart = $d.Paragraphs($para2).Range.Start
$rangeEnd = $d.Paragraphs($para3).Range.End
$rng2 = $d.Range($rangeStart, $rangeEnd)

$xml2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:right="90"/></w:pPr><w:r><w:t>If</w:t></w:r><w:r><w:t xml:space="preserve"> user2 </w:t></w:r><w:r><w:t>accepts</w:t></w:r><w:r><w:t xml:space="preserve"> it sends</w:t></w:r><w:r><w:t xml:space="preserve"> the server connects betwee</w:t></w:r><w:r><w:t>n them by</w:t></w:r><w:r><w:t xml:space="preserve"> sending each user (</w:t></w:r><w:r><w:t>acc,</w:t></w:r><w:r><w:t>user</w:t></w:r><w:r><w:t>uid</w:t></w:r><w:r><w:t>,userip</w:t></w:r><w:r><w:t>:userport</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:left="360" w:right="90" w:firstLine="0"/></w:pPr><w:r><w:t>If user2 rejects it doesn’t send anything</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:ind w:right="90"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng2.InsertXML($xml2)

# The engine's InsertXML round-trip silently drops an explicit
# w:firstLine="0" (it treats it as "unset"/default on re-serialisation).
# Force it back onto the "If user2 rejects ..." paragraph by explicitly
# assigning the ParagraphFormat property (first to a non-zero value so
# the later zero assignment isn't itself treated as a no-op), and
# restore the left indent the same way since the numbered-list
# paragraph's format getter otherwise reports the list's own indent
# instead of the direct override after the XML injection.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*If user2 rejects*") {
        $p.Format.FirstLineIndent = 5
        $p.Format.FirstLineIndent = 0
        $p.Format.LeftIndent = 18
    }
}
